$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 'Purpose
id="Par1">The need for prolonged invasive mechanical ventilation in COVID-19 patients is placing the otorhinolaryngologist in front of an increasing request for tracheostomy.

 Nowadays, there is uncertainty regarding the timing of tracheostomy, the prognosis of these patients and the safety of healthcare workers.

 The aim of this study is to evaluate the efficacy and safety of tracheostomy placement in patients with COVID-19.
Methods
id="Par2">A retrospective cohort study on 23 COVID 19 patients, to analyse the timing of tracheostomy, the risk factors associated with in-hospital death and the infection of the involved health care workers.

 Early tracheostomy was defined as ≤ 10 days and late ones &gt; 10 days.


Results
id="Par3">The mortality rate of COVID-19 patients admitted to ICU that underwent tracheostomy was 18%.

 The overall mortality of patients admitted to ICU was 53%.

 The univariate analysis revealed that early tracheostomy, SOFA score &gt; 6, and D-dimer level &gt; 4 were significantly associated with a greater risk of death.

 At the multivariate analysis SOFA score &gt; 6 and D-dimer level &gt; 4 resulted as significant factors for a higher risk of death.

 No health care workers associated with tracheostomy are confirmed to be infected by SARS-CoV2.
Conclusion
id="Par4">We suggest to wait at least 14 days to perform tracheostomy.

 In patients with SOFA score &gt; 6 and D dimer &gt; 4, tracheostomy should not be performed or should be postponed.

 Optimized procedures and enhanced personal protective equipment can make the tracheostomy safe and beneficial in COVID-19 patients.


'
$ws.Range("E6").Value = '[T.%Volo%tizianavolo@icloud.com%1,       P.%Stritoni%NULL%1,       I.%Battel%NULL%1,       B.%Zennaro%NULL%1,       F.%Lazzari%NULL%1,       M.%Bellin%NULL%1,       L.%Michieletto%NULL%1,       G.%Spinato%NULL%1,       C.%Busatto%NULL%1,       D.%Politi%NULL%1,       R.%Spinato%NULL%1]'
$ws.Range("D10").Value = 'id="Par1">Tracheostomy is performed for critically ill patients with prolonged endotracheal intubation including COVID-19 to mitigate the risk of airway complications.

 We analyzed the difference in mortality of COVID-19 patients undergoing early tracheostomy within 14 days of intubation compared to later tracheostomy than 14 days after intubation.

 The mortality was not significantly higher in early tracheostomy compared to later tracheostomy'
$ws.Range("E10").Value = '[Toshiki%Kuno%kuno-toshiki@hotmail.co.jp%1,       Yoshihisa%Miyamoto%NULL%1,       Matsuo%So%NULL%1,       Masao%Iwagami%NULL%1,       Mai%Takahashi%NULL%1,       Natalia N.%Egorova%NULL%1]'
